$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Notes text (column C) per feedback from the 2nd schema.org workshop ---

# Row 4: url / Dataset-Doi note - tweak wording and the doi URL scheme
$ws.Range("C4").Value = "If Dataset-Doi is not available try Dataset-SecondaryDoi. Any value needs to be prepended with 'https://dx.doi.org/'"

# Row 5: sameAs - expand note with landing-page construction logic
$ws.Range("C5").Value = "Theres a possibility we could construct the dataset landing page URL here. The logic would therefore be 'https://podaac.jpl.nasa.gov/dataset/' + Dataset Shortname e.g. 'UPA-L2P-ATS_NR_2P'"

# Row 12: variableMeasured/description - add a new note that was previously blank
$ws.Range("C12").Value = "Do we want to flag to PO.DAAC Data Engineering that Variable descriptors are lacking???"

# Row 5 grew to two lines worth of text, so it now wraps onto a second line
$ws.Rows.Item(5).RowHeight = 34

# --- Update the view: zoom to 150% and move the selection to C31 ---
$excel.ActiveWindow.Zoom = 150
$ws.Range("C31").Select()
